$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "261.61"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value2 = "1.63%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "27.26"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value2 = "0.60%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "4.719"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value2 = "3.06%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "0.06063"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value2 = "2.81%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "6.639"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value2 = "0.11%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "0.8622"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value2 = "1.27%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.9201"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value2 = "-2.49%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.05123"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value2 = "2.74%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.07097"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value2 = "0.08%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.03044"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value2 = "-0.92%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "0.09098"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value2 = "-0.28%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "0.001534"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value2 = "0.65%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "0.0006100"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value2 = "0.99%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "0.006130"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value2 = "1.76%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value2 = "-1.23%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "3.169"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value2 = "-0.34%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value2 = "2.44%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value2 = "2.21%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "4.092"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value2 = "3.56%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value2 = "-0.24%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "0.001218"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value2 = "-0.35%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value2 = "-8.76%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value2 = "0.07%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value2 = "3.13%"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value2 = "1.38%"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value2 = "1.14%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "0.004132"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value2 = "-33.99%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.01487"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value2 = "5.56%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "0.002182"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value2 = "-0.82%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "0.00005306"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value2 = "-0.76%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value2 = "0.04%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value2 = "-18.42%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "0.1353"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value2 = "-46.24%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value2 = "0.04%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value2 = "0.04%"
